# Experiment 5 header: turn the empty "הרצת ניסויים:" placeholder paragraph
# into a proper sub-heading run ("תרגיל 5:") with the new (smaller, bold)
# run/paragraph formatting + refreshed w14 shadow/outline text effects, and
# relocate the _GoBack bookmark to the following (still-empty) paragraph.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>'
$pkgFooter = '</pkg:xmlData></pkg:part></pkg:package>'
$docOpen = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$docClose = '</w:body></w:document>'

# Shared run/paragraph-mark properties for the new heading text.
$rPr = '<w:rPr><w:rFonts w:cstheme="minorHAnsi" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:rtl/><w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl"><w14:schemeClr w14:val="dk1"><w14:alpha w14:val="60000"/></w14:schemeClr></w14:shadow><w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr>'

# Paragraph that used to hold the stray _GoBack bookmark only; it now gets
# the "תרגיל 5:" run and the updated formatting, and loses the bookmark.
$headingText = "תרגיל 5:"
$para1 = '<w:p w14:paraId="0BD90368" w14:textId="77777777" w:rsidR="008F22EA" w:rsidRDefault="008F22EA" w:rsidP="008F22EA"><w:pPr>' + $rPr + '</w:pPr><w:r>' + $rPr + '<w:t>' + $headingText + '</w:t></w:r></w:p>'

$xml1 = $pkgHeader + $docOpen + $para1 + $docClose + $pkgFooter

$p1 = $d.Paragraphs.Item(4)
$null = $p1.Range.InsertXML($xml1)

# Next (still empty) paragraph now receives the relocated _GoBack bookmark,
# keeping its own original formatting untouched.
$para2 = '<w:p w14:paraId="12138FA0" w14:textId="77777777" w:rsidR="008F22EA" w:rsidRPr="008F22EA" w:rsidRDefault="008F22EA" w:rsidP="008F22EA"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:rtl/><w14:shadow w14:blurRad="38100" w14:dist="25400" w14:dir="5400000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="ctr"><w14:srgbClr w14:val="6E747A"><w14:alpha w14:val="57000"/></w14:srgbClr></w14:shadow><w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$xml2 = $pkgHeader + $docOpen + $para2 + $docClose + $pkgFooter

$p2 = $d.Paragraphs.Item(5)
$null = $p2.Range.InsertXML($xml2)

Write-Output "done"
